$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 99, shifting existing rows 99:158 down to 100:159
$ws.Rows.Item(99).Insert()

# Match the date-column number format used by the other rows in column D
$ws.Range("D99").NumberFormat = $ws.Range("D98").NumberFormat

# Populate the new row 99 with the inserted record's data
$ws.Range("A99").Value = 7
$ws.Range("B99").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C99").Value = "Ñuble"
$ws.Range("D99").Value = 44488
$ws.Range("E99").Value = 16
$ws.Range("F99").Value = 100112043
$ws.Range("G99").Value = "Pepino ensalada"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 60
$ws.Range("K99").Value = 13000
$ws.Range("L99").Value = 14000
$ws.Range("M99").Value = 13500
$ws.Range("N99").Value = "$/caja 60 unidades"
$ws.Range("O99").Value = "Región de Arica y Parinacota"
$ws.Range("P99").Value = 225
$ws.Range("Q99").Value = 60
$ws.Range("R99").Value = "Hortaliza"
